$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 21: finish the "Treeple" section header, mirroring the "YDF"
#     header in row 11 (A21 already contains "Treeple").
$ws.Range("B21").Value = "d"
$ws.Range("C21").Value = "keeping n_attributes = 160"

# --- Row 22: column headers of n_attributes values, mirroring row 12.
$ws.Range("A22").Value = "n"
$ws.Range("B22").Value = 160
$ws.Range("C22").Value = 320
$ws.Range("D22").Value = 640
$ws.Range("E22").Value = 1024
$ws.Range("F22").Value = 2048
$ws.Range("G22").Value = 4096
$ws.Range("H22").Value = 8192

# Copy the bold/italic formatting already used on the equivalent header
# row (row 12) so no new style entries are introduced.
$ws.Range("B12:H12").Copy()
$ws.Range("B22:H22").PasteSpecial(-4122)

# --- Rows 23-27: n values down the left column, matching the formatting
#     used for the same values in rows 13-17.
$ws.Range("A23").Value = 500
$ws.Range("A24").Value = 1000
$ws.Range("A25").Value = 2000
$ws.Range("A26").Value = 4000
$ws.Range("A27").Value = 8000

$ws.Range("A13:A17").Copy()
$ws.Range("A23:A27").PasteSpecial(-4122)

# Note that Treeple couldn't be run for n=500 because the kernel kept
# restarting.
$ws.Range("B23").Value = "Kernel keeps restarting"

$excel.CutCopyMode = 0
$ws.Range("B24").Select()
